$d = $word.ActiveDocument

# The first paragraph currently reads:
#   "This is a Microsoft word document."
# We need to:
#   1. Append two trailing spaces to that run's text.
#   2. Insert a new run after it, colored dark red (C00000), containing:
#      "(This is a change – Version for branch alternate)"

$para1 = $d.Paragraphs(1)
$r = $para1.Range
# Range of the paragraph includes the trailing paragraph mark; back it off
# by one character so we only touch the visible text.
$r.End = $r.End - 1

# Step 1: add two trailing spaces after the existing sentence.
$r.InsertAfter("  ")

# Step 2: insert the new colored run right after the spaces we just added,
# and before the paragraph mark.
$insertPoint = $d.Range($r.End, $r.End)
$insertPoint.InsertAfter("(This is a change " + [char]0x2013 + " Version for branch alternate)")
# Word/VBA color "longs" are packed as 0x00BBGGRR (red in the low byte), the
# same layout the RGB() function produces, so build it that way rather than
# using the OOXML RRGGBB hex order directly.
$insertPoint.Font.Color = 0x000000C0
